$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Fixtures")
$ws2 = $wb.Worksheets.Item("Fixture Types")

# --- New header cells for Fixtures sheet (columns E/F) ---
$ws1.Range("E1").Value = "Universe "
$ws1.Range("F1").Value = "Address"
$ws1.Range("E1:F1").Font.Bold = $true

# --- Fill in Universe/Address data for rows 2-65 ---
$data = @(
    @(2,7,1),
    @(3,7,133),
    @(4,7,142),
    @(5,7,223),
    @(6,7,45),
    @(7,7,169),
    @(8,7,262),
    @(9,7,89),
    @(10,7,196),
    @(11,7,301),
    @(12,8,223),
    @(13,8,142),
    @(14,8,1),
    @(15,8,262),
    @(16,8,169),
    @(17,8,45),
    @(18,8,301),
    @(19,8,196),
    @(20,8,133),
    @(21,8,89),
    @(22,4,1),
    @(23,4,118),
    @(24,4,419),
    @(25,4,199),
    @(26,4,40),
    @(27,4,145),
    @(28,4,428),
    @(29,4,243),
    @(30,4,79),
    @(31,4,172),
    @(32,4,437),
    @(33,4,287),
    @(34,5,287),
    @(35,5,419),
    @(36,5,118),
    @(37,5,1),
    @(38,5,331),
    @(39,5,428),
    @(40,5,145),
    @(41,5,40),
    @(42,5,375),
    @(43,5,437),
    @(44,5,172),
    @(45,5,79),
    @(46,1,1),
    @(47,1,177),
    @(48,1,285),
    @(49,1,45),
    @(50,1,204),
    @(51,1,324),
    @(52,1,89),
    @(53,1,231),
    @(54,1,133),
    @(55,1,258),
    @(56,2,177),
    @(57,2,1),
    @(58,2,204),
    @(59,2,45),
    @(60,2,285),
    @(61,2,231),
    @(62,2,89),
    @(63,2,324),
    @(64,2,258),
    @(65,2,133)
)
foreach ($row in $data) {
    $r = $row[0]
    $e = $row[1]
    $f = $row[2]
    $ws1.Cells.Item($r, 5).Value = $e
    $ws1.Cells.Item($r, 6).Value = $f
}

# --- Widen column D to match new layout ---
$ws1.Columns.Item(4).ColumnWidth = 11.44

# --- Selections / active sheet to match final UI state ---
$null = $ws2.Range("B32").Select()
$ws1.Activate()
$null = $ws1.Range("G3").Select()
